$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.50317888919334
$ws.Range("C2").Value = 11.07942455786873
$ws.Range("D2").Value = 4.944133667519684
$ws.Range("E2").Value = 11.70524948033952
$ws.Range("F2").Value = 24.80316302425875
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("L2").Value = 9.752315673379815
$ws.Range("O2").Value = 22.16467789400161

$ws.Range("B3").Value = 16.84879125733282
$ws.Range("C3").Value = 10.85001874444987
$ws.Range("D3").Value = 4.913707733208355
$ws.Range("E3").Value = 11.76323444131394
$ws.Range("F3").Value = 24.8541947975608
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("L3").Value = 9.719777783250438
$ws.Range("O3").Value = 22.26503839062571

$ws.Range("B4").Value = 16.43507335361609
$ws.Range("C4").Value = 10.70643155393493
$ws.Range("D4").Value = 4.894950157265365
$ws.Range("E4").Value = 11.80116870401728
$ws.Range("F4").Value = 24.89489621898626
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("L4").Value = 9.701404343723134
$ws.Range("O4").Value = 22.33383275546873

$ws.Range("B5").Value = 16.26372241909267
$ws.Range("C5").Value = 10.64729319014809
$ws.Range("D5").Value = 4.887290831353396
$ws.Range("E5").Value = 11.81721351895135
$ws.Range("F5").Value = 24.91382617250817
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("L5").Value = 9.694325670249009
$ws.Range("O5").Value = 22.36366128648234

$ws.Range("B6").Value = 16.23511070378867
$ws.Range("C6").Value = 10.63743738191168
$ws.Range("D6").Value = 4.886018177948373
$ws.Range("E6").Value = 11.81991316183914
$ws.Range("F6").Value = 24.91711066439569
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("L6").Value = 9.693175083761435
$ws.Range("O6").Value = 22.36872238322331

$ws.Range("B7").Value = 16.43277328385415
$ws.Range("C7").Value = 10.70563644253084
$ws.Range("D7").Value = 4.894846918753007
$ws.Range("E7").Value = 11.80138271599648
$ws.Range("F7").Value = 24.89514204168978
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("L7").Value = 9.701307217238291
$ws.Range("O7").Value = 22.33422778110372

$ws.Range("B8").Value = 17.28015647762716
$ws.Range("C8").Value = 11.00092820694904
$ws.Range("D8").Value = 4.933660746441283
$ws.Range("E8").Value = 11.72475888665896
$ws.Range("F8").Value = 24.81880784429977
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("L8").Value = 9.740766662241798
$ws.Range("O8").Value = 22.19778796858003

$ws.Range("B9").Value = 18.83800791761074
$ws.Range("C9").Value = 11.55577663214025
$ws.Range("D9").Value = 5.009013080883248
$ws.Range("E9").Value = 11.59299367616743
$ws.Range("F9").Value = 24.7439058262276
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("L9").Value = 9.830641276665935
$ws.Range("O9").Value = 21.9875485657095

$ws.Range("B10").Value = 19.90849999955306
$ws.Range("C10").Value = 11.94546786629311
$ws.Range("D10").Value = 5.063698408051994
$ws.Range("E10").Value = 11.50745517668962
$ws.Range("F10").Value = 24.73499735289403
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("L10").Value = 9.903944264812464
$ws.Range("O10").Value = 21.86856891161246

$ws.Range("B11").Value = 20.37758168278794
$ws.Range("C11").Value = 12.11823920890677
$ws.Range("D11").Value = 5.088384418048396
$ws.Range("E11").Value = 11.47098805373884
$ws.Range("F11").Value = 24.74103849321733
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("L11").Value = 9.938789744123891
$ws.Range("O11").Value = 21.82225644401392

$ws.Range("B12").Value = 20.55251141460301
$ws.Range("C12").Value = 12.18296948771371
$ws.Range("D12").Value = 5.097701144789268
$ws.Range("E12").Value = 11.45753055120785
$ws.Range("F12").Value = 24.7447812331354
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("L12").Value = 9.952193066900533
$ws.Range("O12").Value = 21.80585097232749

$ws.Range("B13").Value = 20.51495917433762
$ws.Range("C13").Value = 12.16906026678991
$ws.Range("D13").Value = 5.095696076470498
$ws.Range("E13").Value = 11.46041321420474
$ws.Range("F13").Value = 24.74391041298678
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("L13").Value = 9.949297289707184
$ws.Range("O13").Value = 21.8093337013146

$ws.Range("B14").Value = 20.39202807111113
$ws.Range("C14").Value = 12.12357878465577
$ws.Range("D14").Value = 5.089151561002133
$ws.Range("E14").Value = 11.46987384465653
$ws.Range("F14").Value = 24.74131723411847
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("L14").Value = 9.939888317279522
$ws.Range("O14").Value = 21.82088401945856

$ws.Range("B15").Value = 20.31637390677396
$ws.Range("C15").Value = 12.09562827868085
$ws.Range("D15").Value = 5.085138660448555
$ws.Range("E15").Value = 11.47571457739528
$ws.Range("F15").Value = 24.73991840994902
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("L15").Value = 9.934151921195255
$ws.Range("O15").Value = 21.82810660143353

$ws.Range("B16").Value = 19.87747300273855
$ws.Range("C16").Value = 11.93408216381536
$ws.Range("D16").Value = 5.062080947922722
$ws.Range("E16").Value = 11.50988760467902
$ws.Range("F16").Value = 24.73480601242
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("L16").Value = 9.901696568549456
$ws.Range("O16").Value = 21.8717534898097

$ws.Range("B17").Value = 19.60354211103066
$ws.Range("C17").Value = 11.83379132265316
$ws.Range("D17").Value = 5.047884169919006
$ws.Range("E17").Value = 11.53147798287355
$ws.Range("F17").Value = 24.73425810902031
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("L17").Value = 9.882165072762431
$ws.Range("O17").Value = 21.90053658827511

$ws.Range("B18").Value = 19.44430687535157
$ws.Range("C18").Value = 11.77568561968917
$ws.Range("D18").Value = 5.039700849704976
$ws.Range("E18").Value = 11.54412629373983
$ws.Range("F18").Value = 24.73489292762898
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("L18").Value = 9.871072692562558
$ws.Range("O18").Value = 21.91782666785976

$ws.Range("B19").Value = 19.39010867105897
$ws.Range("C19").Value = 11.75594120355339
$ws.Range("D19").Value = 5.036927196375706
$ws.Range("E19").Value = 11.54844830661454
$ws.Range("F19").Value = 24.73527088215578
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("L19").Value = 9.86734155123778
$ws.Range("O19").Value = 21.92380672344672

$ws.Range("B20").Value = 19.63287711771599
$ws.Range("C20").Value = 11.84451138526832
$ws.Range("D20").Value = 5.049397302446549
$ws.Range("E20").Value = 11.52915583442139
$ws.Range("F20").Value = 24.73421808502873
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("L20").Value = 9.88422963039592
$ws.Range("O20").Value = 21.89739646469339

$ws.Range("B21").Value = 20.42821016875106
$ws.Range("C21").Value = 12.13695700583482
$ws.Range("D21").Value = 5.091074725201562
$ws.Range("E21").Value = 11.4670854781151
$ws.Range("F21").Value = 24.74203940287827
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("L21").Value = 9.942646374503138
$ws.Range("O21").Value = 21.81746061996514

$ws.Range("B22").Value = 20.93221310807196
$ws.Range("C22").Value = 12.32401965935733
$ws.Range("D22").Value = 5.118128921544658
$ws.Range("E22").Value = 11.42856964173144
$ws.Range("F22").Value = 24.75563286880459
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("L22").Value = 9.982034329774613
$ws.Range("O22").Value = 21.77182016282529

$ws.Range("B23").Value = 20.66470044453159
$ws.Range("C23").Value = 12.22456769586515
$ws.Range("D23").Value = 5.10370774827059
$ws.Range("E23").Value = 11.44893853586339
$ws.Range("F23").Value = 24.74760096613188
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("L23").Value = 9.960904160279009
$ws.Range("O23").Value = 21.79557245513063

$ws.Range("B24").Value = 19.61962019812058
$ws.Range("C24").Value = 11.83966623417965
$ws.Range("D24").Value = 5.048713281583223
$ws.Range("E24").Value = 11.53020494351291
$ws.Range("F24").Value = 24.73423322151154
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("L24").Value = 9.883295818029572
$ws.Range("O24").Value = 21.89881380297209

$ws.Range("B25").Value = 18.42888426474853
$ws.Range("C25").Value = 11.40861025147862
$ws.Range("D25").Value = 4.988731648774437
$ws.Range("E25").Value = 11.62666038224049
$ws.Range("F25").Value = 24.75610118112119
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("L25").Value = 9.805024285433371
$ws.Range("O25").Value = 22.03823061516892

